$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.983.87"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "3.599.01"
$ws.Range("E3").Value = "  -1.34%  "
$ws.Range("E4").Value = "  -0.28%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.607"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.68%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.665"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.143"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -10.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000250"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -13.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.21%  "
$ws.Range("D14").Value = "4.181.24"
$ws.Range("E14").Value = "  -2.13%  "
$ws.Range("D15").Value = "3.600.41"
$ws.Range("E15").Value = "  -2.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.126"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "66.736.73"
$ws.Range("E17").Value = "  -2.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.08"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "387.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "84.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.95%  "
$ws.Range("E28").Value = "  -11.65%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "65.79"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.77"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.10%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "585.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.35%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.111"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.90"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.60%  "
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -0.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.370"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -6.25%  "
$ws.Range("D40").Value = "0.0₃0726"
$ws.Range("E40").Value = "  -17.83%  "
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("E42").Value = "  -8.47%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0406"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.54%  "
$ws.Range("D44").Value = "2.674.05"
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -12.56%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.129"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.80%  "
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.36%  "
$ws.Range("E48").Value = "  -5.71%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "135.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -9.53%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.57%  "
